$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Type-changing cells (string <-> numeric) ---
$ws.Range("C14").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("C23").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("F31").PasteSpecial(-4122)

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$ws.Range("D33").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("E33").Value = -100
$ws.Range("E15").Copy()
$ws.Range("E33").PasteSpecial(-4122)

$ws.Range("G33").Value = 1
$ws.Range("F14").Copy()
$ws.Range("G33").PasteSpecial(-4122)

$ws.Range("H33").Value = 0
$ws.Range("E15").Copy()
$ws.Range("H33").PasteSpecial(-4122)

# --- Pure value updates ---
$ws.Range("I14").Value = 9
$ws.Range("K14").Value = -52.631578947368
$ws.Range("L14").Value = -18.181818181818
$ws.Range("M14").Value = -43.75
$ws.Range("N14").Value = -59.090909090909
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 60
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 71.428571428571
$ws.Range("M15").Value = 3.448275862068
$ws.Range("N15").Value = -17.808219178082
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 283
$ws.Range("J16").Value = 301
$ws.Range("K16").Value = -5.980066445182
$ws.Range("L16").Value = 6.390977443609
$ws.Range("M16").Value = -30.295566502463
$ws.Range("N16").Value = -78.330781010719
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -11.764705882352
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 72
$ws.Range("H17").Value = -23.611111111111
$ws.Range("I17").Value = 877
$ws.Range("J17").Value = 922
$ws.Range("K17").Value = -4.880694143167
$ws.Range("L17").Value = 16.005291005291
$ws.Range("M17").Value = 92.324561403508
$ws.Range("N17").Value = -21.766280107047
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = -51.219512195122
$ws.Range("I18").Value = 324
$ws.Range("J18").Value = 344
$ws.Range("K18").Value = -5.813953488372
$ws.Range("L18").Value = 2.857142857142
$ws.Range("M18").Value = -45.728643216080
$ws.Range("N18").Value = -89.900249376558
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 52.941176470588
$ws.Range("F19").Value = 82
$ws.Range("G19").Value = 98
$ws.Range("H19").Value = -16.326530612244
$ws.Range("I19").Value = 1363
$ws.Range("J19").Value = 1505
$ws.Range("K19").Value = -9.435215946843
$ws.Range("L19").Value = 3.650190114068
$ws.Range("M19").Value = 47.510822510822
$ws.Range("N19").Value = -13.952020202020
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 322
$ws.Range("J20").Value = 428
$ws.Range("K20").Value = -24.766355140186
$ws.Range("L20").Value = -32.352941176470
$ws.Range("M20").Value = -2.719033232628
$ws.Range("N20").Value = -92.975567190226
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = -15.625
$ws.Range("G21").Value = 277
$ws.Range("H21").Value = -25.270758122743
$ws.Range("I21").Value = 3238
$ws.Range("J21").Value = 3559
$ws.Range("K21").Value = -9.01938746839
$ws.Range("L21").Value = 2.016383112791
$ws.Range("M21").Value = 16.140602582496
$ws.Range("N21").Value = -72.785342074298
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -22.222222222222
$ws.Range("I23").Value = 117
$ws.Range("J23").Value = 128
$ws.Range("K23").Value = -8.59375
$ws.Range("L23").Value = 30
$ws.Range("M23").Value = 72.058823529411
$ws.Range("C24").Value = 74
$ws.Range("D24").Value = 102
$ws.Range("E24").Value = -27.450980392156
$ws.Range("F24").Value = 309
$ws.Range("G24").Value = 404
$ws.Range("H24").Value = -23.514851485148
$ws.Range("I24").Value = 3786
$ws.Range("J24").Value = 4144
$ws.Range("K24").Value = -8.638996138996
$ws.Range("L24").Value = 2.186234817813
$ws.Range("M24").Value = 7.343351290048
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = 19.354838709677
$ws.Range("F25").Value = 171
$ws.Range("G25").Value = 157
$ws.Range("H25").Value = 8.917197452229
$ws.Range("I25").Value = 2042
$ws.Range("J25").Value = 1895
$ws.Range("K25").Value = 7.757255936675
$ws.Range("L25").Value = 52.274422073079
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 36
$ws.Range("E26").Value = -27.777777777777
$ws.Range("F26").Value = 119
$ws.Range("G26").Value = 147
$ws.Range("H26").Value = -19.047619047619
$ws.Range("I26").Value = 1769
$ws.Range("J26").Value = 1670
$ws.Range("K26").Value = 5.928143712574
$ws.Range("L26").Value = 14.055448098001
$ws.Range("M26").Value = -7.187827911857
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -72.727272727272
$ws.Range("I27").Value = 96
$ws.Range("J27").Value = 69
$ws.Range("K27").Value = 39.130434782608
$ws.Range("L27").Value = 29.729729729729
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 14
$ws.Range("H28").Value = -71.428571428571
$ws.Range("J28").Value = 202
$ws.Range("K28").Value = -6.435643564356
$ws.Range("L28").Value = 15.950920245398
$ws.Range("N29").Value = -86.363636363636
$ws.Range("N30").Value = -84.693877551020
$ws.Range("L31").Value = -8.333333333333
$ws.Range("J33").Value = 15
$ws.Range("K33").Value = -33.333333333333
$ws.Range("L33").Value = -16.666666666666

$excel.CutCopyMode = $false